# Insert a new price-record row at row 205 (pushing the existing rows
# 205..265 down to 206..266) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(205).Insert()

$ws.Range("A205").Value = 10
$ws.Range("B205").Value = "Vega Modelo de Temuco"
$ws.Range("C205").Value = "La Araucanía"
$ws.Range("D205").Value = 44809
$ws.Range("E205").Value = 9
$ws.Range("F205").Value = 100112052
$ws.Range("G205").Value = "Albahaca"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 80
$ws.Range("K205").Value = 6000
$ws.Range("L205").Value = 6000
$ws.Range("M205").Value = 6000
$ws.Range("N205").Value = "$/paquete"
$ws.Range("O205").Value = "Región de Arica y Parinacota"
$ws.Range("P205").Value = 6000
$ws.Range("Q205").Value = 1
$ws.Range("R205").Value = "Hortaliza"
